# Apply the "InformacjeOPrzeniesieniach" weekly update:
# - Update the reporting period text on the parameters sheet
# - Replace the table contents on the "Oddzialy" (transfers) sheet with the
#   new set of transferred lessons, shrinking the table from 4 to 3 data rows
# - Narrow column G to fit the new remarks text

$wb = $excel.ActiveWorkbook
$wsParams = $wb.Worksheets.Item(1)
$wsData = $wb.Worksheets.Item(2)

# --- Sheet 1 ("Opis parametrow"): update the period text ---
$wsParams.Range("A2").Value = "Okres: 26.01.2026 (pon.) - 31.01.2026 (sob.)"

# --- Sheet 2 ("Oddzialy"): rewrite the data rows ---

# Row 2 (column G stays blank, as it already is)
$wsData.Range("A2").Value = "26.01.2026, 8, 14:05-14:50, sala: 19"
$wsData.Range("B2").Value = "26.01.2026, 3, 09:40-10:25, sala: 19"
$wsData.Range("C2").Value = "Zaleska Magdalena"
$wsData.Range("D2").Value = "-"
$wsData.Range("E2").Value = "3FA"
$wsData.Range("F2").Value = "Historia"

# Row 3 (column G gets a new remark)
$wsData.Range("A3").Value = "26.01.2026, 10, 15:45-16:30, sala: 31"
$wsData.Range("B3").Value = "26.01.2026, 6, 12:25-13:10, sala: 27"
$wsData.Range("C3").Value = "Kończyńska Małgorzata"
$wsData.Range("D3").Value = "-"
$wsData.Range("E3").Value = "2CB"
$wsData.Range("F3").Value = "Rozwój kompetencji zawodowych - dekoracje w cukiernictwie"
$wsData.Range("G3").Value = "p. Kończyńska, RKZ za lekcję 10"

# Row 4 (column G stays blank, as it already is)
$wsData.Range("A4").Value = "26.01.2026, 7, 13:15-14:00, sala: 40"
$wsData.Range("B4").Value = "26.01.2026, 6, 12:25-13:10, sala: 40"
$wsData.Range("C4").Value = "Socha Dariusz"
$wsData.Range("D4").Value = "-"
$wsData.Range("E4").Value = "2TH"
$wsData.Range("F4").Value = "Obsługa klientów"

# Row 5 no longer exists in the updated table - delete it entirely
$wsData.Range("A5:G5").Delete()

# Narrow column G to fit the shorter remarks text
$wsData.Columns.Item(7).ColumnWidth = 27.3
